# Corrected Calibration and Ingest Sheets for Coastal Gliders
# - Changed FLORT cal value for angular resolution (CC_angular_resolution) to 1.076
# - Changed FLORT cal value for scattering angle (CC_scattering_angle) to 124
# - Left "Asset_Cal_Info" as the active/selected sheet (where the edits were made)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# CC_scattering_angle (row 4, column F)
$ws.Range("F4").Value = 124

# CC_angular_resolution (row 6, column F)
$ws.Range("F6").Value = 1.076

# Make Asset_Cal_Info the active sheet, and update its selected cell
$ws.Activate()
$ws.Range("D25").Select()
